# Regenerate orders with updated distance/size codes.
# Distances: D51 -> D55, D64 -> D69, D80 -> D86
# Size:      S30 -> S31
# These tokens show up embedded inside Condition / Filename_Left /
# Filename_Right / Distance / Size column values (e.g. "Face03_D51_S30",
# "Face03_D51_S30_l.png", "D51", "S30"), so a whole-workbook text
# substitution over the used range reproduces the rename everywhere it
# appears, matching the shared-string table edits in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange

# xlPart = 2 (match any part of cell text), LookAt:=2
# Use positional args in the COM-interop-friendly order: FindWhat, ReplaceWhat, LookAt
$null = $ur.Replace("D51", "D55", 2)
$null = $ur.Replace("D64", "D69", 2)
$null = $ur.Replace("D80", "D86", 2)
$null = $ur.Replace("S30", "S31", 2)
